# Carnet de bord - Stage semaine 2 - Julien LAY
# Fill in the "MERCREDI matin" / "MERCREDI après-midi" block (rows 16-18)
# with the three new activities, resize the affected rows, and move the
# active selection to the last-edited cell (C17:D17), matching the
# author's final save state.
#
# Cell values are entered in the same left-to-right / top-to-bottom order
# the author used (row 16, then row 17 skipping C17, then row 18, then
# finally C17) so the shared-string table comes out in the same order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 16 : "Questionnement sur le réseau de l'entreprise" ---
$ws.Range("B16").Value = "Questionnement sur le réseau de l'entreprise"
$ws.Range("C16").Value = "Questionner un employé de l'entreprise sur le réseau matériel présent dans l'entreprise"
$ws.Range("E16").Value = "Tout le matériel réseau de l'entreprise : 2 NAS, 2 Switch, Une Freebox"
$ws.Range("F16").Value = "30 min"
$ws.Range("G16").Value = "Savoir poser des questions et se référer au cours de réseau de l'IUT"
$ws.Range("I16").Value = "Oberserver, questionner, faire des corrélations avec mes connaissances"
$ws.Range("H16").Value = 3
$ws.Range("J16").Value = 3

# --- Row 17 : "Maintenances de PC" (C17 filled in later, see below) ---
$ws.Range("B17").Value = "Maintenances de PC"
$ws.Range("E17").Value = "Atelier de l'entreprise avec le matériel informatique nécessaire"
$ws.Range("F17").Value = "1h30"
$ws.Range("G17").Value = "Commencer à savoir faire des manipulations simples que l'on m'a montré"
$ws.Range("I17").Value = "Oberserver, questionner, faire des corrélations avec mes connaissances du stage"
$ws.Range("H17").Value = 3
$ws.Range("J17").Value = 4

# --- Row 18 : "Gérer la clientèle" ---
$ws.Range("B18").Value = "Gérer la clientèle"
$ws.Range("C18").Value = "Gérer les demandes de la clientèle ainsi que certaines factures et devis"
$ws.Range("E18").Value = "Poste informatique de l'accueil, imprimante"
$ws.Range("F18").Value = "1h30"
$ws.Range("G18").Value = "Savoir communiquer avec le client et être autonome (ou poser des questions à M. SEGATO)"
$ws.Range("H18").Value = "2 ou 3"
$ws.Range("I18").Value = "Oberserver, questionner, faire des corrélations avec mes connaissances du stage"
$ws.Range("J18").Value = 3

# C17 was typed last
$ws.Range("C17").Value = "Rajouter Windows 10, Office, ninite,com, etc. sur un ordi testé (ne fonctionnait plus avant)"

# Rows grew taller to fit the new wrapped text
$ws.Rows.Item(17).RowHeight = 31.8
$ws.Rows.Item(18).RowHeight = 32.4

# Leave the selection where the author last left it
$ws.Range("C17:D17").Select()
